$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 7: gtnAddButton / companyAdd (Click on id) ---
$ws.Range("C7").Value = "Click"
$ws.Range("D7").Value = "id"
$ws.Range("E7").Value = "gtnAddButton"
$ws.Range("A7").Value = "companyAdd"

# --- Row 8: companyInformation (Click on xpath) ---
$ws.Range("C8").Value = "Click"
$ws.Range("D8").Value = "xpath"
$ws.Range("A8").Value = "companyInformation"
$ws.Range("E8").Value = "//div[contains(text(),'Company Information')]"

# --- Row 9: companyId (Entertext on id) ---
$ws.Range("C9").Value = "Entertext"
$ws.Range("D9").Value = "id"
$ws.Range("E9").Value = "companyInformationTabCompanyId"
$ws.Range("A9").Value = "companyId"

# --- Row 10: companyNo (Entertext on id) ---
$ws.Range("C10").Value = "Entertext"
$ws.Range("D10").Value = "id"
$ws.Range("E10").Value = "companyInformationTabCompanyNo"
$ws.Range("A10").Value = "companyNo"

# --- Row 11: companyName (Entertext on id) ---
$ws.Range("C11").Value = "Entertext"
$ws.Range("D11").Value = "id"
$ws.Range("E11").Value = "companyInformationTabCompanyName"
$ws.Range("A11").Value = "companyName"

# Apply the small monospace "code" style (Consolas 9pt, dark grey) to the
# PropertyName column cells that hold raw element locators, like E4 already
# does for the existing rows.
$codeRange = $ws.Range("E7,E9,E10,E11")
$codeRange.Font.Name = "Consolas"
$codeRange.Font.Family = 3
$codeRange.Font.Size = 9
$codeRange.Font.Color = 2236962

# Move the selection to the last filled cell, and set up the print page.
$ws.Range("A11").Select() | Out-Null
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
